# ---------------------------------------------------------------------------
# Commit: "In pprojector3, debugging KT -> PT"
#
# This script reproduces (on the "thermistor" and "scratch" sheets) the
# edits captured in the target OOXML diff:
#   * thermistor: a new header row (PN / part number) is inserted at the
#     top, the Steinhart-Hart coefficients a/b/c and the temperature are
#     refreshed for the new part, and new Rtot / Current / Power helper
#     cells are added in columns D:F.
#   * scratch: the scratch-pad calculation block is reshuffled - a new
#     calculation block is placed above the previous one, and the
#     previous block is partially cleared out / reused further down.
#   * the workbook-level defined names are re-pointed to the
#     (now one-row-lower) thermistor cells, and "scratch" becomes the
#     active tab.
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

# --- helper style constants -------------------------------------------------
$YELLOW = 65535

# ===========================================================================
# 1. "thermistor" sheet
# ===========================================================================
$th = $wb.Worksheets.Item("thermistor")

# wipe the old layout (content + formatting) so nothing stale is left behind
$th.Range("A1:G20").Clear()

# --- row 1 : new part-number header ----------------------------------------
$th.Range("A1").Value = "PN"
$th.Range("B1").Value = "QTRL2Z-103C3-12"

# --- row 2 : a ---------------------------------------------------------------
$th.Range("A2").Value = "a"
$th.Range("B2").Value = 0.0011127534437700001
$th.Range("B2").NumberFormat = "0.00E+00"
$th.Range("B2").Interior.Color = $YELLOW
$th.Range("D2").Value = "Tmax"
$th.Range("E2").Value = 37

# --- row 3 : b ---------------------------------------------------------------
$th.Range("A3").Value = "b"
$th.Range("B3").Value = 0.000236732362183
$th.Range("B3").NumberFormat = "0.00E+00"
$th.Range("B3").Interior.Color = $YELLOW
$th.Range("D3").Value = "Tmin"
$th.Range("E3").Value = 10

# --- row 4 : c ---------------------------------------------------------------
$th.Range("A4").Value = "c"
$th.Range("B4").Value = 0.000000078005993000000005
$th.Range("B4").NumberFormat = "0.00E+00"
$th.Range("B4").Interior.Color = $YELLOW
$th.Range("D4").Value = "Trange"
$th.Range("E4").Formula = "=E2-E3"

# --- row 5 : T ---------------------------------------------------------------
$th.Range("A5").Value = "T"
$th.Range("B5").Formula = "=37+273"
$th.Range("C5").Value = "K"

# --- row 6 : y = --------------------------------------------------------------
$th.Range("A6").Value = "y ="
$th.Range("B6").Formula = "=(a - 1/T) / SH_c"

# --- row 7 : x = --------------------------------------------------------------
$th.Range("A7").Value = "x ="
$th.Range("B7").Formula = "=SQRT((b/(3*SH_c))^3 + y^2/4)"

# --- row 8 : R = + Rtot ------------------------------------------------------
$th.Range("A8").Value = "R ="
$th.Range("B8").Formula = "= EXP( (x - y/2)^(1/3) - (x+y/2)^(1/3) )"
$th.Range("D8").Value = "Rtot ="
$th.Range("E8").Formula = "=B8+Rref"
$th.Range("F8").Value = "Ohm"

# --- row 9 : Vref + Current --------------------------------------------------
$th.Range("A9").Value = "Vref"
$th.Range("B9").Value = 5
$th.Range("B9").Interior.Color = $YELLOW
$th.Range("C9").Value = "V"
$th.Range("D9").Value = "Current ="
$th.Range("E9").Formula = "=Vref/E8"
$th.Range("E9").NumberFormat = "0.00E+00"
$th.Range("F9").Value = "Amp"

# --- row 10 : Rref + Power ---------------------------------------------------
$th.Range("A10").Value = "Rref"
$th.Range("B10").Value = 10000
$th.Range("B10").Interior.Color = $YELLOW
$th.Range("C10").Value = "Ohm"
$th.Range("D10").Value = "Power ="
$th.Range("E10").Formula = "=E9^2 * B8"
$th.Range("E10").NumberFormat = "0.00E+00"
$th.Range("F10").Value = "Watt"

# --- row 11 : Vadc ------------------------------------------------------------
$th.Range("A11").Value = "Vadc"
$th.Range("B11").Formula = "= B8* Vref / (Rref + B8)"
$th.Range("C11").Value = "V"

# row 12 intentionally left blank

# --- row 13 : measured Vadc swing -------------------------------------------
$th.Range("B13").Formula = "=3.34-1.89"
$th.Range("C13").Value = "V"
$th.Range("D13").Value = "operating range"

# --- row 14 : ADC --------------------------------------------------------------
$th.Range("A14").Value = "ADC"
$th.Range("B14").Value = 12
$th.Range("C14").Value = "bits"

# --- row 15 : effective bits --------------------------------------------------
$th.Range("B15").Formula = "=B14 + LOG(B13/Vref) / LOG(2)"
$th.Range("C15").Value = "bits"
$th.Range("D15").Value = "effectively"

# --- row 16 : final resolution ------------------------------------------------
$th.Range("B16").Formula = "=E4/2^B15"

# --- re-point the workbook defined names to the shifted cells --------------
$wb.Names.Item("a").RefersTo    = "=thermistor!`$B`$2"
$wb.Names.Item("b").RefersTo    = "=thermistor!`$B`$3"
$wb.Names.Item("SH_c").RefersTo = "=thermistor!`$B`$4"
$wb.Names.Item("T").RefersTo    = "=thermistor!`$B`$5"
$wb.Names.Item("y").RefersTo    = "=thermistor!`$B`$6"
$wb.Names.Item("x").RefersTo    = "=thermistor!`$B`$7"
$wb.Names.Item("Vref").RefersTo = "=thermistor!`$B`$9"
$wb.Names.Item("Rref").RefersTo = "=thermistor!`$B`$10"

# ===========================================================================
# 2. "scratch" sheet
# ===========================================================================
$sc = $wb.Worksheets.Item("scratch")

$sc.Range("A1:G20").Clear()

# --- new calculation block (rows 1-3) ---------------------------------------
$sc.Range("A1").Value = 45.57
$sc.Range("B1").Value = 10.37
$sc.Range("C1").Value = 2618
$sc.Range("C1").NumberFormat = "0.00E+00"
$sc.Range("D1").Value = 31.36
$sc.Range("E1").Value = 51002
$sc.Range("E1").Font.Name = "Arial"
$sc.Range("E1").Font.Color = 0

$sc.Range("A2").Value = 48384
$sc.Range("C2").NumberFormat = "0.00E+00"
$sc.Range("D2").NumberFormat = "0.00E+00"

$sc.Range("A3").Formula = "=A1/A2"
$sc.Range("B3").Formula = "=B1/A2"
$sc.Range("C3").Formula = "=C1/A2"
$sc.Range("D3").Formula = "=D1/A2"
$sc.Range("E3").Formula = "=E1/A2"
$sc.Range("A3:E3").NumberFormat = "0.00E+00"

$sc.Range("C4").NumberFormat = "0.00E+00"

# --- previous calculation block, now moved down to rows 5-7 -----------------
$sc.Range("A5").Value = 32.450000000000003
$sc.Range("B5").Value = 31.36
$sc.Range("C5").Value = 2184
$sc.Range("C5").NumberFormat = "0.00E+00"

$sc.Range("A6").Value = 2184

$sc.Range("A7").Formula = "=A5/A6"
$sc.Range("B7").Formula = "=B5/A6"
$sc.Range("C7").Formula = "=C5/A6"
$sc.Range("A7:C7").NumberFormat = "0.00E+00"

# --- small helper values (rows 9-11) ----------------------------------------
$sc.Range("A9").Value = 10.37
$sc.Range("B9").Value = 1

$sc.Range("A10").Value = 2618
$sc.Range("A10").Font.Name = "Arial"
$sc.Range("A10").Font.Color = 0
$sc.Range("A10").NumberFormat = "0.00E+00"
$sc.Range("B10").Font.Name = "Arial"
$sc.Range("B10").Font.Color = 0
$sc.Range("B10").NumberFormat = "0.00E+00"
$sc.Range("C10").NumberFormat = "0.00E+00"
$sc.Range("D10").Font.Name = "Arial"
$sc.Range("D10").Font.Color = 0
$sc.Range("D10").NumberFormat = "0.00E+00"
$sc.Range("E10").NumberFormat = "0.00E+00"
$sc.Range("F10").NumberFormat = "0.00E+00"
$sc.Range("G10").NumberFormat = "0.00E+00"

$sc.Range("A11").Formula = "=A9/A10"
$sc.Range("A11").Font.Name = "Arial"
$sc.Range("A11").Font.Color = 0
$sc.Range("A11").NumberFormat = "0.00E+00"
$sc.Range("B11").Formula = "=B9/A10"
$sc.Range("B11").Font.Name = "Arial"
$sc.Range("B11").Font.Color = 0
$sc.Range("B11").NumberFormat = "0.00E+00"
$sc.Range("C11").NumberFormat = "0.00E+00"
$sc.Range("D11").Font.Name = "Arial"
$sc.Range("D11").Font.Color = 0
$sc.Range("D11").NumberFormat = "0.00E+00"
$sc.Range("E11").NumberFormat = "0.00E+00"
$sc.Range("F11").NumberFormat = "0.00E+00"
$sc.Range("G11").NumberFormat = "0.00E+00"

# --- trailing, still-empty marker cells (rows 13-16) ------------------------
$sc.Range("A13").NumberFormat = "0.00E+00"
$sc.Range("A14").NumberFormat = "0.00E+00"
$sc.Range("A15").NumberFormat = "0.00E+00"
$sc.Range("A16").NumberFormat = "0.00E+00"

# column widths for the new layout
$sc.Columns("A:A").ColumnWidth = 9.166666666666666
$sc.Columns("B:B").ColumnWidth = 8.6

# ===========================================================================
# 3. Selections / active sheet
# ===========================================================================
$th.Range("E9").Select()

$sc.Activate()
$sc.Range("B12").Select()
